# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 146
    3  = 417
    4  = 12074
    6  = 129
    8  = 86
    10 = 184
    11 = 434
    12 = 54
    16 = 350
    17 = 2209
    18 = 85
    19 = 926
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
